$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 180330
$ws.Range("C4").Value = 170282
$ws.Range("C7").Value = 5.57
$ws.Range("C8").Value = 65.23999999999999
